# Auto-generated Excel COM-interop script
# Applies numeric corrections to the Leve profit-tracking sheets
# (Exodus_Profits workbook) per the scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 845.8200000000001
$ws.Range("I15").Value = 845.8200000000001
$ws.Range("K15").Value = 2537.46
$ws.Range("M15").Value = -2368.46
# Row 18
$ws.Range("H18").Value = 1306.7142
$ws.Range("I18").Value = 1132.6842
$ws.Range("J18").Value = 2960
$ws.Range("K18").Value = 1132.6842
$ws.Range("L18").Value = 2960
$ws.Range("M18").Value = -848.6841999999999
$ws.Range("N18").Value = -3528
# Row 40
$ws.Range("H40").Value = 9828.714
$ws.Range("I40").Value = 5549.75
$ws.Range("J40").Value = 15534
$ws.Range("K40").Value = 5549.75
$ws.Range("L40").Value = 15534
$ws.Range("M40").Value = -5374.75
$ws.Range("N40").Value = -15884
# Row 76
$ws.Range("H76").Value = 33335884
$ws.Range("J76").Value = 3329.8
$ws.Range("L76").Value = 3329.8
$ws.Range("N76").Value = -3959.8
# Row 79
$ws.Range("H79").Value = 33335884
$ws.Range("J79").Value = 3329.8
$ws.Range("L79").Value = 3329.8
$ws.Range("N79").Value = -5513.8
# Row 112
$ws.Range("H112").Value = 1198.3478
$ws.Range("I112").Value = 520.6
$ws.Range("J112").Value = 1386.6111
$ws.Range("K112").Value = 1561.8
$ws.Range("L112").Value = 4159.8333
$ws.Range("M112").Value = -453.8000000000002
$ws.Range("N112").Value = -6375.8333
# Row 137
$ws.Range("H137").Value = 455217.5
$ws.Range("I137").Value = 1863.5264
$ws.Range("K137").Value = 5590.5792
$ws.Range("M137").Value = -3040.5792

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5896.902
$ws.Range("I32").Value = 2796.976
$ws.Range("K32").Value = 2796.976
$ws.Range("M32").Value = -2509.976
# Row 61
$ws.Range("H61").Value = 113027.336
$ws.Range("I61").Value = 1624.5
$ws.Range("J61").Value = 335833
$ws.Range("K61").Value = 1624.5
$ws.Range("L61").Value = 335833
$ws.Range("M61").Value = -1412.5
$ws.Range("N61").Value = -336257
# Row 74
$ws.Range("H74").Value = 129359.5
$ws.Range("I74").Value = 203455.4
$ws.Range("J74").Value = 5866.3335
$ws.Range("K74").Value = 203455.4
$ws.Range("L74").Value = 5866.3335
$ws.Range("M74").Value = -202581.4
$ws.Range("N74").Value = -7614.3335
# Row 76
$ws.Range("H76").Value = 116666
$ws.Range("J76").Value = 116666
$ws.Range("L76").Value = 116666
$ws.Range("N76").Value = -117342
# Row 77
$ws.Range("H77").Value = 129359.5
$ws.Range("I77").Value = 203455.4
$ws.Range("J77").Value = 5866.3335
$ws.Range("K77").Value = 1017277
$ws.Range("L77").Value = 29331.6675
$ws.Range("M77").Value = -1012909
$ws.Range("N77").Value = -38067.6675
# Row 79
$ws.Range("H79").Value = 116666
$ws.Range("J79").Value = 116666
$ws.Range("L79").Value = 116666
$ws.Range("N79").Value = -119006
# Row 110
$ws.Range("H110").Value = 1198.7333
$ws.Range("I110").Value = 1107.25
$ws.Range("J110").Value = 1564.6666
$ws.Range("K110").Value = 1107.25
$ws.Range("L110").Value = 1564.6666
$ws.Range("M110").Value = 937.75
$ws.Range("N110").Value = -5654.6666
# Row 122
$ws.Range("H122").Value = 2928.75
$ws.Range("I122").Value = 2864.5
$ws.Range("K122").Value = 8593.5
$ws.Range("M122").Value = -6143.5
# Row 136
$ws.Range("H136").Value = 113027.336
$ws.Range("I136").Value = 1624.5
$ws.Range("J136").Value = 335833
$ws.Range("K136").Value = 4873.5
$ws.Range("L136").Value = 1007499
$ws.Range("M136").Value = -2323.5
$ws.Range("N136").Value = -1012599

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 181511.16
$ws.Range("I20").Value = 214251.9
$ws.Range("J20").Value = 1437
$ws.Range("K20").Value = 214251.9
$ws.Range("L20").Value = 1437
$ws.Range("M20").Value = -214004.9
$ws.Range("N20").Value = -1931
# Row 94
$ws.Range("H94").Value = 2397.1292
$ws.Range("I94").Value = 2363.4443
$ws.Range("K94").Value = 2363.4443
$ws.Range("M94").Value = -1912.4443
# Row 99
$ws.Range("H99").Value = 2225800.2
$ws.Range("I99").Value = 113136
$ws.Range("K99").Value = 113136
$ws.Range("M99").Value = -111638
# Row 107
$ws.Range("H107").Value = 2418.0715
$ws.Range("I107").Value = 1922.1111
$ws.Range("K107").Value = 1922.1111
$ws.Range("M107").Value = -2.111100000000079

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 1000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
# Row 27
$ws.Range("H27").Value = 1000
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
# Row 58
$ws.Range("H58").Value = 1678.1052
$ws.Range("I58").Value = 1390.3846
$ws.Range("K58").Value = 1390.3846
$ws.Range("M58").Value = -1187.3846
# Row 99
$ws.Range("H99").Value = 8376464
$ws.Range("I99").Value = 13891049
$ws.Range("K99").Value = 13891049
$ws.Range("M99").Value = -13889551
# Row 126
$ws.Range("H126").Value = 8376464
$ws.Range("I126").Value = 13891049
$ws.Range("K126").Value = 41673147
$ws.Range("M126").Value = -41670677
# Row 132
$ws.Range("H132").Value = 2382.65
$ws.Range("I132").Value = 2078.8572
$ws.Range("K132").Value = 6236.571599999999
$ws.Range("M132").Value = -3706.571599999999
# Row 134
$ws.Range("H134").Value = 30811.361
$ws.Range("I134").Value = 3213.0356
$ws.Range("K134").Value = 9639.106800000001
$ws.Range("M134").Value = -7104.106800000001
# Row 136
$ws.Range("H136").Value = 1678.1052
$ws.Range("I136").Value = 1390.3846
$ws.Range("K136").Value = 4171.1538
$ws.Range("M136").Value = -1621.1538

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 128.55
$ws.Range("I97").Value = 112.69231
$ws.Range("J97").Value = 158
$ws.Range("K97").Value = 338.07693
$ws.Range("L97").Value = 474
$ws.Range("M97").Value = 157.92307
$ws.Range("N97").Value = -1466
# Row 117
$ws.Range("H117").Value = 404.66666
$ws.Range("J117").Value = 99.5
$ws.Range("L117").Value = 298.5
$ws.Range("N117").Value = -7182.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7000
$ws.Range("I70").Value = 7000
$ws.Range("K70").Value = 7000
$ws.Range("M70").Value = -6730
# Row 73
$ws.Range("H73").Value = 7000
$ws.Range("I73").Value = 7000
$ws.Range("K73").Value = 7000
$ws.Range("M73").Value = -6064
# Row 97
$ws.Range("H97").Value = 721.2857
$ws.Range("I97").Value = 721.2857
$ws.Range("K97").Value = 721.2857
$ws.Range("M97").Value = -225.2857
# Row 102
$ws.Range("H102").Value = 2220
$ws.Range("I102").Value = 2117.5
$ws.Range("K102").Value = 2117.5
$ws.Range("M102").Value = -495.5
# Row 113
$ws.Range("H113").Value = 2154822
$ws.Range("J113").Value = 6669065.5
$ws.Range("L113").Value = 6669065.5
$ws.Range("N113").Value = -6673405.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 8171
$ws.Range("I22").Value = 1199
$ws.Range("K22").Value = 1199
$ws.Range("M22").Value = -904
# Row 27
$ws.Range("H27").Value = 8171
$ws.Range("I27").Value = 1199
$ws.Range("K27").Value = 1199
$ws.Range("M27").Value = -1092
# Row 46
$ws.Range("H46").Value = 7625.778
$ws.Range("I46").Value = 13574.125
$ws.Range("J46").Value = 2867.1
$ws.Range("K46").Value = 13574.125
$ws.Range("L46").Value = 2867.1
$ws.Range("M46").Value = -13386.125
$ws.Range("N46").Value = -3243.1
# Row 55
$ws.Range("H55").Value = 8484.733
$ws.Range("J55").Value = 16842.715
$ws.Range("L55").Value = 16842.715
$ws.Range("N55").Value = -17188.715
# Row 61
$ws.Range("H61").Value = 1036
$ws.Range("I61").Value = 771
$ws.Range("K61").Value = 771
$ws.Range("M61").Value = -569
# Row 74
$ws.Range("H74").Value = 29833.334
$ws.Range("I74").Value = 19750
$ws.Range("K74").Value = 19750
$ws.Range("M74").Value = -18752
# Row 77
$ws.Range("H77").Value = 29833.334
$ws.Range("I77").Value = 19750
$ws.Range("K77").Value = 59250
$ws.Range("M77").Value = -54258
# Row 113
$ws.Range("H113").Value = 1036
$ws.Range("I113").Value = 771
$ws.Range("K113").Value = 771
$ws.Range("M113").Value = 1399
# Row 132
$ws.Range("H132").Value = 2925
$ws.Range("I132").Value = 2925
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8775
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6245
$ws.Range("N132").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 114235.25
$ws.Range("J46").Value = 114235.25
$ws.Range("L46").Value = 114235.25
$ws.Range("N46").Value = -114697.25
# Row 122
$ws.Range("H122").Value = 2462.3333
$ws.Range("I122").Value = 2278.4285
$ws.Range("K122").Value = 6835.2855
$ws.Range("M122").Value = -4385.2855
# Row 132
$ws.Range("H132").Value = 2498.4443
$ws.Range("I132").Value = 2331.8333
$ws.Range("K132").Value = 6995.499899999999
$ws.Range("M132").Value = -4465.499899999999
# Row 134
$ws.Range("H134").Value = 114235.25
$ws.Range("J134").Value = 114235.25
$ws.Range("L134").Value = 342705.75
$ws.Range("N134").Value = -347775.75

Write-Host "Applied 225 value updates and 3 cell clears."
